$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 609 (shifts existing rows 609.. down to 613..)
$ws.Rows("609:612").Insert()

# Row 609: Murcott / Especial, new week (2023-10-19), volumen changed 400 -> 700
$ws.Range("A609").Value = 2
$ws.Range("B609").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C609").Value = "Coquimbo"
$ws.Range("D609").Value = 45218
$ws.Range("E609").Value = 4
$ws.Range("F609").Value = "Fruta"
$ws.Range("G609").Value = 100102
$ws.Range("H609").Value = "Cítricos"
$ws.Range("I609").Value = 100102004
$ws.Range("J609").Value = "Mandarina"
$ws.Range("K609").Value = "Murcott"
$ws.Range("L609").Value = "Especial"
$ws.Range("M609").Value = 700
$ws.Range("N609").Value = 5500
$ws.Range("O609").Value = 6000
$ws.Range("P609").Value = 5750
$ws.Range("Q609").Value = "$/bandeja 10 kilos"
$ws.Range("R609").Value = "Provincia de Limarí"
$ws.Range("S609").Value = 575
$ws.Range("T609").Value = 10

# Row 610: Murcott / Primera, new week (2023-10-19), volumen changed 400 -> 1100
$ws.Range("A610").Value = 2
$ws.Range("B610").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C610").Value = "Coquimbo"
$ws.Range("D610").Value = 45218
$ws.Range("E610").Value = 4
$ws.Range("F610").Value = "Fruta"
$ws.Range("G610").Value = 100102
$ws.Range("H610").Value = "Cítricos"
$ws.Range("I610").Value = 100102004
$ws.Range("J610").Value = "Mandarina"
$ws.Range("K610").Value = "Murcott"
$ws.Range("L610").Value = "Primera"
$ws.Range("M610").Value = 1100
$ws.Range("N610").Value = 4500
$ws.Range("O610").Value = 5000
$ws.Range("P610").Value = 4750
$ws.Range("Q610").Value = "$/bandeja 10 kilos"
$ws.Range("R610").Value = "Provincia de Limarí"
$ws.Range("S610").Value = 475
$ws.Range("T610").Value = 10

# Row 611: Murcott / Segunda, new week (2023-10-19), volumen changed 360 -> 800
$ws.Range("A611").Value = 2
$ws.Range("B611").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C611").Value = "Coquimbo"
$ws.Range("D611").Value = 45218
$ws.Range("E611").Value = 4
$ws.Range("F611").Value = "Fruta"
$ws.Range("G611").Value = 100102
$ws.Range("H611").Value = "Cítricos"
$ws.Range("I611").Value = 100102004
$ws.Range("J611").Value = "Mandarina"
$ws.Range("K611").Value = "Murcott"
$ws.Range("L611").Value = "Segunda"
$ws.Range("M611").Value = 800
$ws.Range("N611").Value = 3500
$ws.Range("O611").Value = 4000
$ws.Range("P611").Value = 3750
$ws.Range("Q611").Value = "$/bandeja 10 kilos"
$ws.Range("R611").Value = "Provincia de Limarí"
$ws.Range("S611").Value = 375
$ws.Range("T611").Value = 10

# Row 612: Murcott / Tercera, new week (2023-10-19), volumen + prices changed
$ws.Range("A612").Value = 2
$ws.Range("B612").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C612").Value = "Coquimbo"
$ws.Range("D612").Value = 45218
$ws.Range("E612").Value = 4
$ws.Range("F612").Value = "Fruta"
$ws.Range("G612").Value = 100102
$ws.Range("H612").Value = "Cítricos"
$ws.Range("I612").Value = 100102004
$ws.Range("J612").Value = "Mandarina"
$ws.Range("K612").Value = "Murcott"
$ws.Range("L612").Value = "Tercera"
$ws.Range("M612").Value = 500
$ws.Range("N612").Value = 2500
$ws.Range("O612").Value = 3000
$ws.Range("P612").Value = 2750
$ws.Range("Q612").Value = "$/bandeja 10 kilos"
$ws.Range("R612").Value = "Provincia de Limarí"
$ws.Range("S612").Value = 275
$ws.Range("T612").Value = 10
